$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Rows 16 and 18 (periods 2410 and 2502) had their "Periodo Mora" (E) and
# "Valor Mora" (F) values swapped between them; row 17 (2412) is unchanged.
$ws.Range("E16").Value = "2502"
$ws.Range("F16").Value = 56940

$ws.Range("E18").Value = "2410"
$ws.Range("F18").Value = 52000
